$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value2 = 3083.4
$ws.Range("I2").Value2 = 311.8
$ws.Range("J2").Value2 = 5855
$ws.Range("K2").Value2 = 311.8
$ws.Range("L2").Value2 = 5855
$ws.Range("M2").Value2 = -198.8
$ws.Range("N2").Value2 = -6081

# Row 32
$ws.Range("H32").Value2 = 2857.0908
$ws.Range("I32").Value2 = 2699.5
$ws.Range("J32").Value2 = 2947.1428
$ws.Range("K32").Value2 = 2699.5
$ws.Range("L32").Value2 = 2947.1428
$ws.Range("M32").Value2 = -2373.5
$ws.Range("N32").Value2 = -3599.1428

# Row 40
$ws.Range("H40").Value2 = 1666.4
$ws.Range("I40").Value2 = 1687.68
$ws.Range("K40").Value2 = 1687.68
$ws.Range("M40").Value2 = -1512.68

# Row 62
$ws.Range("H62").Value2 = 2998
$ws.Range("I62").Value2 = 2998
$ws.Range("K62").Value2 = 2998
$ws.Range("M62").Value2 = -2374

# Row 65
$ws.Range("H65").Value2 = 2998
$ws.Range("I65").Value2 = 2998
$ws.Range("K65").Value2 = 14990
$ws.Range("M65").Value2 = -11870

# Row 69
$ws.Range("H69").Value2 = 31909.47
$ws.Range("I69").Value2 = 113006.5
$ws.Range("J69").Value2 = 21096.533
$ws.Range("K69").Value2 = 339019.5
$ws.Range("L69").Value2 = 63289.599
$ws.Range("M69").Value2 = -338145.5
$ws.Range("N69").Value2 = -65037.599

# Row 72
$ws.Range("H72").Value2 = 31909.47
$ws.Range("I72").Value2 = 113006.5
$ws.Range("J72").Value2 = 21096.533
$ws.Range("K72").Value2 = 1017058.5
$ws.Range("L72").Value2 = 189868.797
$ws.Range("M72").Value2 = -1012690.5
$ws.Range("N72").Value2 = -198604.797

# Row 96
$ws.Range("H96").Value2 = 5262.4
$ws.Range("I96").Value2 = 7462.3335
$ws.Range("J96").Value2 = 1962.5
$ws.Range("K96").Value2 = 22387.0005
$ws.Range("L96").Value2 = 5887.5
$ws.Range("M96").Value2 = -21014.0005
$ws.Range("N96").Value2 = -8633.5

# Row 97
$ws.Range("H97").Value2 = 2385.8333
$ws.Range("J97").Value2 = 2385.8333
$ws.Range("L97").Value2 = 7157.499899999999
$ws.Range("N97").Value2 = -8149.499899999999

# Row 99
$ws.Range("H99").Value2 = 181
$ws.Range("I99").Value2 = 0
$ws.Range("K99").Value2 = 0
$ws.Range("M99").ClearContents()

# Row 101
$ws.Range("H101").Value2 = 10003227
$ws.Range("I101").Value2 = 16669999
$ws.Range("J101").Value2 = 3069.75
$ws.Range("K101").Value2 = 50009997
$ws.Range("L101").Value2 = 9209.25
$ws.Range("M101").Value2 = -50008375
$ws.Range("N101").Value2 = -12453.25

# Row 107
$ws.Range("H107").Value2 = 1225.5
$ws.Range("I107").Value2 = 450
$ws.Range("K107").Value2 = 450
$ws.Range("M107").Value2 = 1470

# Row 115
$ws.Range("H115").Value2 = 2320.5
$ws.Range("I115").Value2 = 685
$ws.Range("J115").Value2 = 2865.6667
$ws.Range("K115").Value2 = 2055
$ws.Range("L115").Value2 = 8597.000100000001
$ws.Range("M115").Value2 = -488
$ws.Range("N115").Value2 = -11731.0001

# Row 116
$ws.Range("H116").Value2 = 3709.4075
$ws.Range("I116").Value2 = 3004.5386
$ws.Range("J116").Value2 = 4363.9287
$ws.Range("K116").Value2 = 3004.5386
$ws.Range("L116").Value2 = 4363.9287
$ws.Range("M116").Value2 = 437.4614000000001
$ws.Range("N116").Value2 = -11247.9287

# Row 132
$ws.Range("H132").Value2 = 2125.2307
$ws.Range("I132").Value2 = 2125.2307
$ws.Range("K132").Value2 = 6375.6921
$ws.Range("M132").Value2 = -3845.6921

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value2 = 9573.5
$ws.Range("J46").Value2 = 9514.666999999999
$ws.Range("L46").Value2 = 9514.666999999999
$ws.Range("N46").Value2 = -10152.667

# Row 61
$ws.Range("H61").Value2 = 5194
$ws.Range("I61").Value2 = 5247.6665
$ws.Range("K61").Value2 = 5247.6665
$ws.Range("M61").Value2 = -5035.6665

# Row 88
$ws.Range("H88").Value2 = 2021.4286
$ws.Range("J88").Value2 = 1573.5
$ws.Range("L88").Value2 = 1573.5
$ws.Range("N88").Value2 = -2385.5

# Row 91
$ws.Range("H91").Value2 = 2021.4286
$ws.Range("J91").Value2 = 1573.5
$ws.Range("L91").Value2 = 1573.5
$ws.Range("N91").Value2 = -4381.5

# Row 101
$ws.Range("H101").Value2 = 20000
$ws.Range("I101").Value2 = 20000
$ws.Range("J101").Value2 = 0
$ws.Range("K101").Value2 = 20000
$ws.Range("L101").Value2 = 0
$ws.Range("M101").Value2 = -16755
$ws.Range("N101").ClearContents()

# Row 132
$ws.Range("H132").Value2 = 2189.4167
$ws.Range("I132").Value2 = 2184.875
$ws.Range("J132").Value2 = 2198.5
$ws.Range("K132").Value2 = 6554.625
$ws.Range("L132").Value2 = 6595.5
$ws.Range("M132").Value2 = -4024.625
$ws.Range("N132").Value2 = -11655.5

# Row 136
$ws.Range("H136").Value2 = 5194
$ws.Range("I136").Value2 = 5247.6665
$ws.Range("K136").Value2 = 15742.9995
$ws.Range("M136").Value2 = -13192.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value2 = 1751.75
$ws.Range("I11").Value2 = 1335.6666
$ws.Range("J11").Value2 = 3000
$ws.Range("K11").Value2 = 1335.6666
$ws.Range("L11").Value2 = 3000
$ws.Range("M11").Value2 = -1195.6666
$ws.Range("N11").Value2 = -3280

# Row 94
$ws.Range("H94").Value2 = 1671.6316
$ws.Range("I94").Value2 = 1625.8889
$ws.Range("K94").Value2 = 1625.8889
$ws.Range("M94").Value2 = -1174.8889

# Row 134
$ws.Range("H134").Value2 = 3685.9473
$ws.Range("I134").Value2 = 3835.1667
$ws.Range("K134").Value2 = 11505.5001
$ws.Range("M134").Value2 = -8970.500100000001

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value2 = 686
$ws.Range("I17").Value2 = 686
$ws.Range("K17").Value2 = 686
$ws.Range("M17").Value2 = -512

# Row 132
$ws.Range("H132").Value2 = 7534.9443
$ws.Range("I132").Value2 = 7986.125
$ws.Range("J132").Value2 = 3925.5
$ws.Range("K132").Value2 = 23958.375
$ws.Range("L132").Value2 = 11776.5
$ws.Range("M132").Value2 = -21428.375
$ws.Range("N132").Value2 = -16836.5

# Row 134
$ws.Range("H134").Value2 = 0
$ws.Range("I134").Value2 = 0
$ws.Range("K134").Value2 = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value2 = 131.33333
$ws.Range("I14").Value2 = 131.33333
$ws.Range("K14").Value2 = 393.99999
$ws.Range("M14").Value2 = -220.99999

# Row 40
$ws.Range("H40").Value2 = 99.5
$ws.Range("J40").Value2 = 311.66666
$ws.Range("L40").Value2 = 1246.66664
$ws.Range("N40").Value2 = -1384.66664

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value2 = 2624
$ws.Range("I102").Value2 = 2246.5715
$ws.Range("K102").Value2 = 2246.5715
$ws.Range("M102").Value2 = -624.5715

# Row 107
$ws.Range("H107").Value2 = 1201
$ws.Range("I107").Value2 = 900
$ws.Range("J107").Value2 = 1502
$ws.Range("K107").Value2 = 900
$ws.Range("L107").Value2 = 1502
$ws.Range("M107").Value2 = 1020
$ws.Range("N107").Value2 = -5342

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value2 = 1686.9166
$ws.Range("I46").Value2 = 1830.5
$ws.Range("K46").Value2 = 1830.5
$ws.Range("M46").Value2 = -1642.5

# Row 136
$ws.Range("H136").Value2 = 2529.3333
$ws.Range("I136").Value2 = 2372.25
$ws.Range("K136").Value2 = 7116.75
$ws.Range("M136").Value2 = -4566.75
